# Daily burndown chart update
# Record today's "Actual" value (C4) to match the already-recorded
# "Planned" value (B4), and move the selection forward to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 27

$ws.Range("C9").Select()
